# Commit: "last commit & move rapport"
# The results table (methods x candidates, rows 4-12, columns B:AQ) is
# overwritten so that every cell now references the candidate "Montebourg",
# with the single exception of cell H7 which references "Philipot".
# (The shared-strings table correspondingly shrinks to just the method/label
# strings plus these two candidate names.)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B4:AQ12").Value = "Montebourg"
$ws.Range("H7").Value = "Philipot"
